# Update column F ("dSF") values on the active worksheet to match the
# repulled data / mean calculation described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 4
    5  = 1
    6  = -4
    8  = -1
    9  = -3
    10 = -2
    11 = -2
    12 = -4
    13 = 7
    14 = -1
    15 = -1
    16 = 2
    17 = 8
    18 = 0
    20 = 4
    21 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
